$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 144.5
$ws.Range("I9").Value = 144.5
$ws.Range("K9").Value = 144.5
$ws.Range("M9").Value = 24.5
$ws.Range("H96").Value = 779.6667
$ws.Range("I96").Value = 655
$ws.Range("J96").Value = 1029
$ws.Range("K96").Value = 1965
$ws.Range("L96").Value = 3087
$ws.Range("M96").Value = -592
$ws.Range("N96").Value = -5833
$ws.Range("H111").Value = 500000
$ws.Range("I111").Value = 500000
$ws.Range("K111").Value = 1500000
$ws.Range("M111").Value = -1496933
$ws.Range("H125").Value = 56108.055
$ws.Range("I125").Value = 167171.83
$ws.Range("J125").Value = 576.1667
$ws.Range("K125").Value = 1504546.47
$ws.Range("L125").Value = 5185.5003
$ws.Range("M125").Value = -1502086.47
$ws.Range("N125").Value = -10105.5003

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 40000
$ws.Range("J44").Value = 40000
$ws.Range("L44").Value = 40000
$ws.Range("N44").Value = -40976
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H74").Value = 1167.7241
$ws.Range("I74").Value = 962.5
$ws.Range("J74").Value = 1623.7778
$ws.Range("K74").Value = 962.5
$ws.Range("L74").Value = 1623.7778
$ws.Range("M74").Value = -88.5
$ws.Range("N74").Value = -3371.7778
$ws.Range("H77").Value = 1167.7241
$ws.Range("I77").Value = 962.5
$ws.Range("J77").Value = 1623.7778
$ws.Range("K77").Value = 4812.5
$ws.Range("L77").Value = 8118.889
$ws.Range("M77").Value = -444.5
$ws.Range("N77").Value = -16854.889
$ws.Range("H110").Value = 884.36365
$ws.Range("I110").Value = 937.9
$ws.Range("J110").Value = 349
$ws.Range("K110").Value = 937.9
$ws.Range("L110").Value = 349
$ws.Range("M110").Value = 1107.1
$ws.Range("N110").Value = -4439
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H131").Value = 30237.666
$ws.Range("J131").Value = 30237.666
$ws.Range("L131").Value = 30237.666
$ws.Range("N131").Value = -40317.666

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H129").Value = 37249.5
$ws.Range("J129").Value = 37249.5
$ws.Range("L129").Value = 37249.5
$ws.Range("N129").Value = -47249.5
$ws.Range("H130").Value = 57530
$ws.Range("J130").Value = 57530
$ws.Range("L130").Value = 57530
$ws.Range("N130").Value = -67570
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H94").Value = 2002.1578
$ws.Range("I94").Value = 1050
$ws.Range("J94").Value = 2256.0667
$ws.Range("K94").Value = 1050
$ws.Range("L94").Value = 2256.0667
$ws.Range("M94").Value = -599
$ws.Range("N94").Value = -3158.0667
$ws.Range("H99").Value = 3189.9443
$ws.Range("I99").Value = 3189.9443
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3189.9443
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1691.9443
$ws.Range("N99").ClearContents()
$ws.Range("H110").Value = 20750
$ws.Range("J110").Value = 20750
$ws.Range("L110").Value = 20750
$ws.Range("N110").Value = -28930
$ws.Range("H122").Value = 4390.1113
$ws.Range("I122").Value = 4390.1113
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 13170.3339
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -10720.3339
$ws.Range("N122").ClearContents()
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H124").Value = 24700
$ws.Range("J124").Value = 24700
$ws.Range("L124").Value = 24700
$ws.Range("N124").Value = -29610
$ws.Range("H125").Value = 30000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 30000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 30000
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -34920
$ws.Range("H126").Value = 3189.9443
$ws.Range("I126").Value = 3189.9443
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9569.832900000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7099.832900000001
$ws.Range("N126").ClearContents()
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H131").Value = 10000
$ws.Range("I131").Value = 10000
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 10000
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -4960
$ws.Range("N131").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 278.58823
$ws.Range("I17").Value = 248
$ws.Range("K17").Value = 744
$ws.Range("M17").Value = -575
$ws.Range("H34").Value = 1555.7142
$ws.Range("I34").Value = 630
$ws.Range("J34").Value = 2250
$ws.Range("K34").Value = 1890
$ws.Range("L34").Value = 6750
$ws.Range("M34").Value = -1806
$ws.Range("N34").Value = -6918
$ws.Range("H39").Value = 2516
$ws.Range("J39").Value = 2516
$ws.Range("L39").Value = 7548
$ws.Range("N39").Value = -8136
$ws.Range("H55").Value = 10637.846
$ws.Range("J55").Value = 11482.667
$ws.Range("L55").Value = 34448.001
$ws.Range("N55").Value = -34802.001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 35000
$ws.Range("J92").Value = 35000
$ws.Range("L92").Value = 35000
$ws.Range("N92").Value = -38744
$ws.Range("H126").Value = 3768
$ws.Range("I126").Value = 1812
$ws.Range("J126").Value = 4257
$ws.Range("K126").Value = 5436
$ws.Range("L126").Value = 12771
$ws.Range("M126").Value = -2966
$ws.Range("N126").Value = -17711

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9928
$ws.Range("J2").Value = 9928
$ws.Range("L2").Value = 9928
$ws.Range("N2").Value = -10152
$ws.Range("H57").Value = 4100
$ws.Range("I57").Value = 4100
$ws.Range("K57").Value = 4100
$ws.Range("M57").Value = -3534
$ws.Range("H64").Value = 19500
$ws.Range("J64").Value = 19500
$ws.Range("L64").Value = 19500
$ws.Range("N64").Value = -19950
$ws.Range("H67").Value = 19500
$ws.Range("J67").Value = 19500
$ws.Range("L67").Value = 19500
$ws.Range("N67").Value = -21060
$ws.Range("H97").Value = 30249.25
$ws.Range("J97").Value = 30249.25
$ws.Range("L97").Value = 30249.25
$ws.Range("N97").Value = -32231.25
$ws.Range("H136").Value = 2352.6453
$ws.Range("I136").Value = 2034.619
$ws.Range("K136").Value = 6103.857
$ws.Range("M136").Value = -3553.857

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 152500
$ws.Range("J95").Value = 152500
$ws.Range("L95").Value = 152500
$ws.Range("N95").Value = -157992
